$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Header rename: STEP NUMBER -> LINE NUMBER
$ws.Range("A1").Value = "LINE NUMBER"

# 2. Insert a new row at row 33 (pushes existing rows 33.. down by one)
$ws.Rows.Item(33).Insert()

# 3. Populate the newly inserted row 33 with the "iterate (after while)" flow-type entry
$ws.Cells.Item(33, 1).Value = 6
$ws.Cells.Item(33, 2).Value = "flow type"
$ws.Cells.Item(33, 3).Value = "iterate  (after while)"
